# Error Calculations and Plots
# Applies the "missing_data.xlsx" imputation edits:
#  - removes the two rows whose ID is "RM 232" and "SC 92" (rows shift up)
#  - fills in / blanks out a number of individual numeric cells elsewhere
#    on the sheet (simulating additional values being marked missing, and
#    some previously-missing values being imputed/filled back in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the "RM 232" row (original row 26) and the "SC 92" row
#    (original row 28, which becomes row 27 once row 26 is gone).
# ---------------------------------------------------------------------
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# ---------------------------------------------------------------------
# 2) Per-cell value edits among the rows that did not shift (rows 2-25).
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 4).Value = -14.2

$ws.Cells.Item(4, 5).ClearContents()

$ws.Cells.Item(5, 4).ClearContents()
$ws.Cells.Item(5, 6).ClearContents()

$ws.Cells.Item(9, 5).Value = -6.8

$ws.Cells.Item(10, 5).Value = -6.1

$ws.Cells.Item(11, 5).Value = -7.9

$ws.Cells.Item(12, 5).Value = -5.3
$ws.Cells.Item(12, 6).ClearContents()

$ws.Cells.Item(13, 6).Value = 17.1

$ws.Cells.Item(14, 6).Value = 17.76

$ws.Cells.Item(15, 5).ClearContents()

$ws.Cells.Item(17, 5).ClearContents()
$ws.Cells.Item(17, 6).Value = 17.78

$ws.Cells.Item(18, 5).ClearContents()

$ws.Cells.Item(20, 5).ClearContents()
$ws.Cells.Item(20, 6).ClearContents()

$ws.Cells.Item(21, 4).Value = -14.3

$ws.Cells.Item(23, 4).ClearContents()
$ws.Cells.Item(23, 6).ClearContents()

$ws.Cells.Item(25, 6).ClearContents()

# ---------------------------------------------------------------------
# 3) Per-cell value edits among the rows that shifted up after the two
#    row deletions above (using the FINAL row numbers, 26-33).
# ---------------------------------------------------------------------
$ws.Cells.Item(27, 6).Value = 17

$ws.Cells.Item(28, 6).Value = 17.44

$ws.Cells.Item(31, 5).Value = -8.1

$ws.Cells.Item(32, 4).Value = -14.7
$ws.Cells.Item(32, 5).Value = -6.4
$ws.Cells.Item(32, 6).Value = 17.39
